$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header style (bold, border, centered) from existing header E1 to F1:H1
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean outlier flag data for rows 2-12, columns F (KNN), G (SVM), H (RF)
$data = @(
    @(0,0,0),
    @(1,1,1),
    @(0,0,1),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,1),
    @(0,0,0),
    @(0,0,0)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 6).Value = [bool]$vals[0]
    $ws.Cells.Item($row, 7).Value = [bool]$vals[1]
    $ws.Cells.Item($row, 8).Value = [bool]$vals[2]
    $row++
}
